# "Generate Report for Handoff"
#
# The localization CI run moved this item from "In Translation" to
# "Ready for handoff" and re-generated the xliff hand-off files, which
# bumped the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps. Saving the regenerated report also re-flowed the "Status"
# columns' auto-fit width on the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ------------------
$wsOverview.Range("E2").Value2 = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value2 = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value2     = "Ready for handoff"
$wsDeDe.Range("C2").Value2     = "Ready for handoff"

# --- Timestamps bumped by the new handoff generation -------------------
# Overview!G2 "Latest HO Xliff Generate Date" and de-de!H2
# "Latest Handoff Datetime" both advance to 13:08:21.
$wsOverview.Range("G2").Value2 = "2016-09-02 13:08:21"
$wsDeDe.Range("H2").Value2     = "2016-09-02 13:08:21"

# zh-cn!H2 "Latest Handoff Datetime" advances to 13:08:15.
$wsZhCn.Range("H2").Value2     = "2016-09-02 13:08:15"

# --- Status column widened (re-fit for the longer "Ready for handoff") -
$wsOverview.Columns.Item(5).ColumnWidth = 16.3   # E:E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 16.3   # F:F (de-de status)
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.3   # C:C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.3   # C:C (Status)
